$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get shuffled between rows (by column letter => column index)
# D=4 (Fecha), L=12 (Calidad), M=13 (Volumen), N=14 (Precio minimo),
# O=15 (Precio maximo), P=16 (Precio promedio ponderado),
# Q=17 (Unidad de comercializacion), S=19 (Precio $/Kg), T=20 (Kg/unidad)
$cols = @(4, 12, 13, 14, 15, 16, 17, 19, 20)

# Mapping: new row number -> row number whose data (in the shuffled columns)
# should be copied into it. Rows 1 (header) and all other columns are
# untouched.
$rowMap = @{
    2  = 4
    3  = 5
    4  = 2
    5  = 3
    6  = 8
    7  = 15
    8  = 13
    9  = 11
    10 = 12
    11 = 9
    12 = 14
    13 = 10
    14 = 6
    15 = 7
}

# Snapshot the current (pre-edit) values of the shuffled columns for every
# data row before writing anything back, since several rows swap with each
# other (e.g. 2<->4, 3<->5, 7<->15, 9<->11) and some form longer cycles
# (6->8->13->10->12->14->6).
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($r in $rowMap.Keys) {
    $src = $rowMap[$r]
    $rowVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c]
    }
}
